$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Header row additions
$ws.Range("F1").Value = "HS Graduate"
$ws.Range("G1").Value = "Age"

# Column F ("HS Graduate") values
$fValues = @("Yes","yes","yes","yes","no","na","yes","no","no","yes","yes","yes","yes","yes")
for ($i = 0; $i -lt $fValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $fValues[$i]
}

# Column G ("Age") values
$gValues = @(25,30,19,23,"na",17,44,15,16,50,48,41,23,35)
for ($i = 0; $i -lt $gValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $gValues[$i]
}

# Center-align the whole used range (A1:G15)
$ws.Range("A1:G15").HorizontalAlignment = -4108

# Bold the new header cell through column F to match the rest of the header row
$ws.Range("F1").Font.Bold = $true

# Autofit columns A-F so column widths best-fit their content (mirrors Excel's bestFit columns)
$ws.Range("A1:F15").EntireColumn.AutoFit()

$ws.Range("G16").Select()
